$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "added new line in main line"
$ws.Range("A3").Select()
